# Actualización automática de tasas-transfi.xlsx

$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the daily conversion text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$check = [string][char]0x2705

$oldText = $wsHoja1.Range("A1").Value()
$newText = $oldText.Replace(
    ($check + " 1000 Bs = 3.28 = 12544.26 pesos"),
    ($check + " 1000 Bs = 3.24 = 12401.39 pesos")
)
$newText = $newText.Replace(
    ($check + " 12544.26 pesos = 3.26 = 969.0 Bs"),
    ($check + " 12401.39 pesos = 3.21 = 970.27 Bs")
)
$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update rate values ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 309
$wsTasas.Range("O10").Value = 3832.03
$wsTasas.Range("N12").Value = 3859.96
$wsTasas.Range("O12").Value = 302
